$d = $word.ActiveDocument

# 1) Bump the header-row heights of the two affected tables (571 -> 637 twips,
#    i.e. 28.55pt -> 31.85pt), keeping the existing "auto" height rule.
#    Re-resolve each table anchor immediately before mutating it (rather than
#    caching both anchors up front) so each edit is applied to the live
#    document state.
$d.Tables.Item(2).Rows.First.Height = 637 / 20
$d.Tables.Item(4).Rows.First.Height = 637 / 20

# 2) Fix the mangled "chi" glyph used in the "chi-squared" column headers:
#    replace every occurrence of the correct "χ" (U+03C7) with the mojibake
#    "Ï‡" (U+00CF U+2021) that now appears in the source.
$d.Content.Find.Execute("χ", $false, $false, $false, $false, $false, $true, 1, $false, "Ï‡", 2)
